$d = $word.ActiveDocument

# --- Locate the three paragraphs we need to touch -------------------------
# "In order to win, a row can have cells that are:"
# "Visible and empty (mineless)"
# "Hidden and have a mine"
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("In order to win")) {
        $targetIdx = $i
        break
    }
}

$pWin = $d.Paragraphs.Item($targetIdx)
$pVisible = $d.Paragraphs.Item($targetIdx + 1)
$pHidden = $d.Paragraphs.Item($targetIdx + 2)

# --- 1) Rewrite "In order to win, a row can have cells that are:" as the --
#        six-run sentence "In order to be complete, all the empty(mineless)
#        cells in a row must be visible" (keeping its ListParagraph /
#        numPr ilvl=1 numId=3 paragraph properties intact). We use
#        Range.InsertXML so each chunk lands in its own <w:r>, matching the
#        target run layout exactly instead of being coalesced into one run.
$winXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>In order to</w:t></w:r><w:r><w:t xml:space="preserve"> be</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>complete</w:t></w:r><w:r><w:t>, a</w:t></w:r><w:r><w:t>ll the empty(mineless) cells in a row must be visible</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $pWin.Range.InsertXML($winXml)

# --- 2) "Visible and empty (mineless)" loses its bullet (numPr) and its ---
#        text; it becomes an empty ListParagraph-styled paragraph indented
#        to where the old bullet text used to sit (ind left=1800).
$emptyXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1800"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $pVisible.Range.InsertXML($emptyXml)

# --- 3) "Hidden and have a mine" paragraph is removed entirely (text + ----
#        its paragraph mark).
$null = $pHidden.Range.Delete()
